# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed from the regenerated
# underlying stats (std/mean -> s_vals) and rewritten into the save_data
# workbook. Only the numeric contents of column G change; everything else
# (headers, other columns, styles, layout) stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, as produced by the regenerated
# s_vals calculation.
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 2
    9  = 3
    10 = 2
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 3
    17 = 1
    18 = 1
    19 = 0
    20 = 3
    21 = 0
    22 = 3
    23 = 1
    24 = 0
    25 = 2
    26 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 0
    36 = 0
    37 = 1
    38 = 0
    41 = 2
    42 = 2
    44 = 2
    45 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
